# Autogenerated on Thu Mar 26 2015 18:06:15 GMT+0000 (Coordinated Universal Time)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Remove the existing hyperlink (ref A49 -> http://www.stat.kg/rus/part/msp.htm)
$ws.Hyperlinks.Delete()

# Insert a new blank row above the old row 48 ("National Statistical Committee").
# This shifts rows 48-54 down to 49-55 and the newly inserted row 48 naturally
# inherits the "source" (italic) formatting already used by the rows around it.
$ws.Rows("48:48").Insert()

# Remember the hyperlink text before it is overwritten (currently sitting at A50
# after the shift, still carrying the old hyperlink formatting).
$urlText = $ws.Range("A50").Text

# The url text moves one row further down, to A51, which already has the correct
# plain "source" (italic) formatting inherited from the row shift above.
$ws.Range("A51").Value = $urlText

# A50 becomes the blank separator line; clear its old hyperlink formatting so it
# matches the plain italic "source" style used by the rest of this block.
$ws.Range("A50").Value = ""
$ws.Range("A50").Font.Italic = $true
$ws.Range("A50").Font.Underline = $false
$ws.Range("A50").Font.ColorIndex = -4105

# Update the long citation text (now located at A55 after the row shift).
$ws.Range("A55").Value = "National Statistical Committee of the Kyrgyz Republic (NSCKR), ""Краткие методологические пояснения"", 2014. Available at http://www.stat.kg/stat.files/tematika/%D0%A4%D0%B8%D0%BD%D0%B0%D0%BD%D1%81%D1%8B/msp/%D0%BF%D0%BE%D1%8F%D1%81%D0%BD%D0%B5%D0%BD%D0%B8%D1%8F%20%D0%BF%D0%BE%20%D0%9C%D0%A1%D0%9F.pdf"
